# Auto-generated edit script: updates PYLL Excel data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1820.6
$ws.Range("D2").Value = 275.7
$ws.Range("E2").Value = 16979.5
$ws.Range("F2").Value = 4418.3
$ws.Range("G2").Value = 9.33
$ws.Range("H2").Value = 0.88
$ws.Range("J2").Value = 1202.9
$ws.Range("K2").Value = 313
$ws.Range("L2").Value = "16979.5 (±4418.3)"
$ws.Range("M2").Value = "9.33 (±0.88)"
$ws.Range("N2").Value = "1202.9 (±313.0)"

# Row 3
$ws.Range("C3").Value = 3774.400000000001
$ws.Range("D3").Value = 339.2
$ws.Range("E3").Value = 42985.4
$ws.Range("F3").Value = 5586.8
$ws.Range("G3").Value = 11.39
$ws.Range("H3").Value = 0.42
$ws.Range("J3").Value = 3052.7
$ws.Range("K3").Value = 396.8
$ws.Range("L3").Value = "42985.4 (±5586.8)"
$ws.Range("M3").Value = "11.39 (±0.42)"
$ws.Range("N3").Value = "3052.7 (±396.8)"

# Row 4
$ws.Range("C4").Value = 5593.6
$ws.Range("D4").Value = 507.4
$ws.Range("E4").Value = 59923.3
$ws.Range("F4").Value = 8541.900000000001
$ws.Range("G4").Value = 10.71
$ws.Range("H4").Value = 0.51
$ws.Range("J4").Value = 2125.2
$ws.Range("K4").Value = 302.9
$ws.Range("L4").Value = "59923.3 (±8541.9)"
$ws.Range("M4").Value = "10.71 (±0.51)"
$ws.Range("N4").Value = "2125.2 (±302.9)"

# Row 11
$ws.Range("C11").Value = 5336.200000000001
$ws.Range("E11").Value = 56598.8
$ws.Range("G11").Value = 10.61
$ws.Range("H11").Value = 0.61
$ws.Range("J11").Value = 983.3
$ws.Range("L11").Value = "56598.8 (±14943.3)"
$ws.Range("M11").Value = "10.61 (±0.61)"
$ws.Range("N11").Value = "983.3 (±259.6)"

# Row 12
$ws.Range("C12").Value = 6663.4
$ws.Range("E12").Value = 71872.70000000001
$ws.Range("J12").Value = 1272.9
$ws.Range("L12").Value = "71872.7 (±18243.5)"
$ws.Range("M12").Value = "10.79 (±0.6)"
$ws.Range("N12").Value = "1272.9 (±323.1)"

# Row 13
$ws.Range("C13").Value = 11959.4
$ws.Range("E13").Value = 129756.6
$ws.Range("J13").Value = 1138
$ws.Range("L13").Value = "129756.6 (±31186.9)"
$ws.Range("M13").Value = "10.85 (±0.62)"
$ws.Range("N13").Value = "1138.0 (±273.5)"

# Row 14
$ws.Range("C14").Value = 6379.4
$ws.Range("E14").Value = 87261.39999999999
$ws.Range("J14").Value = 2455.6
$ws.Range("L14").Value = "87261.4 (±14421.8)"
$ws.Range("M14").Value = "13.68 (±0.24)"
$ws.Range("N14").Value = "2455.6 (±405.8)"

# Row 16
$ws.Range("C16").Value = 15983.8
$ws.Range("E16").Value = 215233.4
$ws.Range("J16").Value = 3114.5
$ws.Range("L16").Value = "215233.4 (±30250.1)"
$ws.Range("M16").Value = "13.47 (±0.11)"
$ws.Range("N16").Value = "3114.5 (±437.7)"

# Row 23
$ws.Range("C23").Value = 6171.4
$ws.Range("E23").Value = 73373.7
$ws.Range("J23").Value = 1365.2
$ws.Range("L23").Value = "73373.7 (±18448.1)"
$ws.Range("M23").Value = "11.89 (±0.29)"
$ws.Range("N23").Value = "1365.2 (±343.2)"

# Row 24
$ws.Range("C24").Value = 9570
$ws.Range("E24").Value = 103511.1
$ws.Range("G24").Value = 10.82
$ws.Range("H24").Value = 0.49
$ws.Range("J24").Value = 1969.8
$ws.Range("L24").Value = "103511.1 (±20823.3)"
$ws.Range("M24").Value = "10.82 (±0.49)"
$ws.Range("N24").Value = "1969.8 (±396.3)"

# Row 25
$ws.Range("C25").Value = 15732
$ws.Range("E25").Value = 182739.9
$ws.Range("G25").Value = 11.62
$ws.Range("H25").Value = 0.36
$ws.Range("J25").Value = 1719.2
$ws.Range("L25").Value = "182739.9 (±37384.4)"
$ws.Range("M25").Value = "11.62 (±0.36)"
$ws.Range("N25").Value = "1719.2 (±351.7)"

# Row 29
$ws.Range("C29").Value = 293.4
$ws.Range("E29").Value = 5501.3
$ws.Range("G29").Value = 18.75
$ws.Range("H29").Value = 1.45
$ws.Range("J29").Value = 797.1
$ws.Range("L29").Value = "5501.3 (±2909.1)"
$ws.Range("M29").Value = "18.75 (±1.45)"
$ws.Range("N29").Value = "797.1 (±421.5)"

# Row 30
$ws.Range("C30").Value = 375.9999999999999
$ws.Range("E30").Value = 5213.7
$ws.Range("G30").Value = 13.87
$ws.Range("H30").Value = 1.17
$ws.Range("J30").Value = 831.4
$ws.Range("L30").Value = "5213.7 (±2819.9)"
$ws.Range("M30").Value = "13.87 (±1.17)"
$ws.Range("N30").Value = "831.4 (±449.7)"

# Row 31
$ws.Range("C31").Value = 627.8
$ws.Range("E31").Value = 10417.5
$ws.Range("G31").Value = 16.59
$ws.Range("H31").Value = 0.95
$ws.Range("J31").Value = 790.9
$ws.Range("L31").Value = "10417.5 (±5399.3)"
$ws.Range("M31").Value = "16.59 (±0.95)"
$ws.Range("N31").Value = "790.9 (±409.9)"

# Row 35
$ws.Range("C35").Value = 13328.6
$ws.Range("E35").Value = 183543.9
$ws.Range("J35").Value = 538.2
$ws.Range("L35").Value = "183543.9 (±49339.3)"
$ws.Range("M35").Value = "13.77 (±0.1)"
$ws.Range("N35").Value = "538.2 (±144.7)"

# Row 36
$ws.Range("C36").Value = 23570.4
$ws.Range("E36").Value = 274447.4
$ws.Range("J36").Value = 850.3
$ws.Range("L36").Value = "274447.4 (±66866.6)"
$ws.Range("M36").Value = "11.64 (±0.46)"
$ws.Range("N36").Value = "850.3 (±207.2)"

# Row 37
$ws.Range("C37").Value = 36878
$ws.Range("E37").Value = 471649.1
$ws.Range("J37").Value = 710.5
$ws.Range("L37").Value = "471649.1 (±116464.7)"
$ws.Range("M37").Value = "12.79 (±0.32)"
$ws.Range("N37").Value = "710.5 (±175.5)"

# Row 51
$ws.Range("C51").Value = 28039.2
$ws.Range("E51").Value = 313670.3
$ws.Range("J51").Value = 1076.9
$ws.Range("L51").Value = "313670.3 (±70253.6)"
$ws.Range("M51").Value = "11.19 (±0.21)"
$ws.Range("N51").Value = "1076.9 (±241.2)"

# Row 52
$ws.Range("C52").Value = 45990.6
$ws.Range("E52").Value = 527596.6
$ws.Range("H52").Value = 0.25
$ws.Range("J52").Value = 1862.6
$ws.Range("L52").Value = "527596.6 (±66277.4)"
$ws.Range("M52").Value = "11.47 (±0.25)"
$ws.Range("N52").Value = "1862.6 (±234.0)"

# Row 53
$ws.Range("C53").Value = 74025.60000000001
$ws.Range("E53").Value = 869892.2
$ws.Range("J53").Value = 1478.2
$ws.Range("L53").Value = "869892.2 (±132785.4)"
$ws.Range("M53").Value = "11.75 (±0.21)"
$ws.Range("N53").Value = "1478.2 (±225.6)"

# Row 54
$ws.Range("C54").Value = 414.8
$ws.Range("D54").Value = 189.2
$ws.Range("E54").Value = 3881.5
$ws.Range("F54").Value = 1922.2
$ws.Range("G54").Value = 9.359999999999999
$ws.Range("H54").Value = 0.25
$ws.Range("J54").Value = 382.5
$ws.Range("K54").Value = 189.4
$ws.Range("L54").Value = "3881.5 (±1922.2)"
$ws.Range("M54").Value = "9.36 (±0.25)"
$ws.Range("N54").Value = "382.5 (±189.4)"

# Row 55
$ws.Range("C55").Value = 459.4
$ws.Range("D55").Value = 211
$ws.Range("E55").Value = 4645.9
$ws.Range("F55").Value = 2719.3
$ws.Range("G55").Value = 10.11
$ws.Range("H55").Value = 0.88
$ws.Range("J55").Value = 528.9
$ws.Range("K55").Value = 309.6
$ws.Range("L55").Value = "4645.9 (±2719.3)"
$ws.Range("M55").Value = "10.11 (±0.88)"
$ws.Range("N55").Value = "528.9 (±309.6)"

# Row 56
$ws.Range("C56").Value = 838.5999999999999
$ws.Range("E56").Value = 8439.4
$ws.Range("G56").Value = 10.06
$ws.Range("H56").Value = 0.48
$ws.Range("J56").Value = 445.8
$ws.Range("L56").Value = "8439.4 (±5006.9)"
$ws.Range("M56").Value = "10.06 (±0.48)"
$ws.Range("N56").Value = "445.8 (±264.5)"

# Row 63
$ws.Range("C63").Value = 160.8
$ws.Range("E63").Value = 2095.8
$ws.Range("G63").Value = 13.03
$ws.Range("H63").Value = 2.15
$ws.Range("J63").Value = 850.1
$ws.Range("L63").Value = "2095.8 (±1419.4)"
$ws.Range("M63").Value = "13.03 (±2.15)"
$ws.Range("N63").Value = "850.1 (±575.7)"

# Row 64
$ws.Range("C64").Value = 283.8
$ws.Range("E64").Value = 3949.7
$ws.Range("G64").Value = 13.92
$ws.Range("H64").Value = 1.11
$ws.Range("J64").Value = 1492.2
$ws.Range("L64").Value = "3949.7 (±1788.5)"
$ws.Range("M64").Value = "13.92 (±1.11)"
$ws.Range("N64").Value = "1492.2 (±675.7)"

# Row 65
$ws.Range("C65").Value = 438.6
$ws.Range("E65").Value = 6007.7
$ws.Range("G65").Value = 13.7
$ws.Range("H65").Value = 1.41
$ws.Range("J65").Value = 1175.2
$ws.Range("L65").Value = "6007.7 (±2785.3)"
$ws.Range("M65").Value = "13.7 (±1.41)"
$ws.Range("N65").Value = "1175.2 (±544.8)"

# Row 66
$ws.Range("C66").Value = 288.4
$ws.Range("E66").Value = 3640.400000000001
$ws.Range("G66").Value = 12.62
$ws.Range("H66").Value = 1.01
$ws.Range("J66").Value = 1163.5
$ws.Range("L66").Value = "3640.4 (±1393.6)"
$ws.Range("M66").Value = "12.62 (±1.01)"
$ws.Range("N66").Value = "1163.5 (±445.4)"

# Row 67
$ws.Range("C67").Value = 521.2
$ws.Range("E67").Value = 5686.599999999999
$ws.Range("G67").Value = 10.91
$ws.Range("H67").Value = 0.9399999999999999
$ws.Range("J67").Value = 1855.8
$ws.Range("L67").Value = "5686.6 (±2367.4)"
$ws.Range("M67").Value = "10.91 (±0.94)"
$ws.Range("N67").Value = "1855.8 (±772.6)"

# Row 68
$ws.Range("C68").Value = 788.6
$ws.Range("E68").Value = 9356.199999999999
$ws.Range("G68").Value = 11.86
$ws.Range("H68").Value = 1.05
$ws.Range("J68").Value = 1510.7
$ws.Range("L68").Value = "9356.2 (±2887.2)"
$ws.Range("M68").Value = "11.86 (±1.05)"
$ws.Range("N68").Value = "1510.7 (±466.2)"

# Row 69
$ws.Range("C69").Value = 5709.4
$ws.Range("D69").Value = 1060.7
$ws.Range("E69").Value = 65278.8
$ws.Range("F69").Value = 18124.5
$ws.Range("H69").Value = 0.89
$ws.Range("J69").Value = 753.4
$ws.Range("K69").Value = 209.2
$ws.Range("L69").Value = "65278.8 (±18124.5)"
$ws.Range("M69").Value = "11.43 (±0.89)"
$ws.Range("N69").Value = "753.4 (±209.2)"

# Row 70
$ws.Range("C70").Value = 9518.799999999999
$ws.Range("D70").Value = 1417
$ws.Range("E70").Value = 104431.1
$ws.Range("F70").Value = 20808.4
$ws.Range("G70").Value = 10.97
$ws.Range("H70").Value = 0.48
$ws.Range("J70").Value = 1213.3
$ws.Range("K70").Value = 241.7
$ws.Range("L70").Value = "104431.1 (±20808.4)"
$ws.Range("M70").Value = "10.97 (±0.48)"
$ws.Range("N70").Value = "1213.3 (±241.7)"

# Row 71
$ws.Range("C71").Value = 15133
$ws.Range("D71").Value = 2378.2
$ws.Range("E71").Value = 169793.3
$ws.Range("F71").Value = 39057.2
$ws.Range("G71").Value = 11.22
$ws.Range("H71").Value = 0.71
$ws.Range("J71").Value = 983.1
$ws.Range("K71").Value = 226.1
$ws.Range("L71").Value = "169793.3 (±39057.2)"
$ws.Range("M71").Value = "11.22 (±0.71)"
$ws.Range("N71").Value = "983.1 (±226.1)"

# Row 72
$ws.Range("C72").Value = 297.6
$ws.Range("D72").Value = 236.9
$ws.Range("E72").Value = 4613.7
$ws.Range("F72").Value = 3964.2
$ws.Range("G72").Value = 15.5
$ws.Range("H72").Value = 0.55
$ws.Range("J72").Value = 175.5
$ws.Range("K72").Value = 150.8
$ws.Range("L72").Value = "4613.7 (±3964.2)"
$ws.Range("M72").Value = "15.5 (±0.55)"
$ws.Range("N72").Value = "175.5 (±150.8)"

# Row 73
$ws.Range("C73").Value = 458
$ws.Range("E73").Value = 6373.5
$ws.Range("F73").Value = 3249.8
$ws.Range("G73").Value = 13.92
$ws.Range("H73").Value = 0.23
$ws.Range("J73").Value = 236.7
$ws.Range("K73").Value = 120.7
$ws.Range("L73").Value = "6373.5 (±3249.8)"
$ws.Range("M73").Value = "13.92 (±0.23)"
$ws.Range("N73").Value = "236.7 (±120.7)"

# Row 74
$ws.Range("C74").Value = 745.5999999999999
$ws.Range("D74").Value = 458.6
$ws.Range("E74").Value = 10587.1
$ws.Range("F74").Value = 6956.299999999999
$ws.Range("G74").Value = 14.2
$ws.Range("H74").Value = 0.37
$ws.Range("J74").Value = 198.9
$ws.Range("K74").Value = 130.7
$ws.Range("L74").Value = "10587.1 (±6956.3)"
$ws.Range("M74").Value = "14.2 (±0.37)"
$ws.Range("N74").Value = "198.9 (±130.7)"

# Row 81
$ws.Range("C81").Value = 13183.4
$ws.Range("E81").Value = 170212.7
$ws.Range("J81").Value = 1740.8
$ws.Range("L81").Value = "170212.7 (±33669.7)"
$ws.Range("M81").Value = "12.91 (±0.05)"
$ws.Range("N81").Value = "1740.8 (±344.3)"

# Row 82
$ws.Range("C82").Value = 21703.4
$ws.Range("E82").Value = 273195.1
$ws.Range("J82").Value = 2901.3
$ws.Range("L82").Value = "273195.1 (±43481.5)"
$ws.Range("M82").Value = "12.59 (±0.11)"
$ws.Range("N82").Value = "2901.3 (±461.8)"

# Row 83
$ws.Range("C83").Value = 34886.8
$ws.Range("E83").Value = 463350.2999999999
$ws.Range("J83").Value = 2414
$ws.Range("L83").Value = "463350.3 (±75790.6)"
$ws.Range("M83").Value = "13.28 (±0.01)"
$ws.Range("N83").Value = "2414.0 (±394.9)"

# Row 87
$ws.Range("C87").Value = 2457.6
$ws.Range("E87").Value = 31991.4
$ws.Range("J87").Value = 1153.9
$ws.Range("L87").Value = "31991.4 (±9693.9)"
$ws.Range("M87").Value = "13.02 (±0.69)"
$ws.Range("N87").Value = "1153.9 (±349.7)"

# Row 88
$ws.Range("C88").Value = 3420
$ws.Range("E88").Value = 38055.2
$ws.Range("G88").Value = 11.13
$ws.Range("J88").Value = 1431.8
$ws.Range("L88").Value = "38055.2 (±7216.4)"
$ws.Range("M88").Value = "11.13 (±0.43)"
$ws.Range("N88").Value = "1431.8 (±271.5)"

# Row 89
$ws.Range("C89").Value = 5834.199999999999
$ws.Range("E89").Value = 71073.8
$ws.Range("J89").Value = 1308.8
$ws.Range("L89").Value = "71073.8 (±11401.1)"
$ws.Range("M89").Value = "12.18 (±0.13)"
$ws.Range("N89").Value = "1308.8 (±210.0)"

# Row 90
$ws.Range("C90").Value = 1168.4
$ws.Range("E90").Value = 12666.4
$ws.Range("G90").Value = 10.84
$ws.Range("H90").Value = 0.76
$ws.Range("J90").Value = 1228.8
$ws.Range("L90").Value = "12666.4 (±3540.4)"
$ws.Range("M90").Value = "10.84 (±0.76)"
$ws.Range("N90").Value = "1228.8 (±343.5)"

# Row 91
$ws.Range("C91").Value = 1514.8
$ws.Range("E91").Value = 16816.5
$ws.Range("G91").Value = 11.1
$ws.Range("H91").Value = 0.92
$ws.Range("J91").Value = 1607
$ws.Range("L91").Value = "16816.5 (±4708.2)"
$ws.Range("M91").Value = "11.1 (±0.92)"
$ws.Range("N91").Value = "1607.0 (±449.9)"

# Row 92
$ws.Range("C92").Value = 2671
$ws.Range("E92").Value = 29858.7
$ws.Range("G92").Value = 11.18
$ws.Range("H92").Value = 0.6
$ws.Range("J92").Value = 1437.4
$ws.Range("L92").Value = "29858.7 (±6247.2)"
$ws.Range("M92").Value = "11.18 (±0.6)"
$ws.Range("N92").Value = "1437.4 (±300.7)"

# Row 93
$ws.Range("C93").Value = 21162.2
$ws.Range("D93").Value = 3315.8
$ws.Range("E93").Value = 277662.6
$ws.Range("F93").Value = 48046.60000000001
$ws.Range("G93").Value = 13.12
$ws.Range("H93").Value = 0.19
$ws.Range("J93").Value = 1170.4
$ws.Range("K93").Value = 202.5
$ws.Range("L93").Value = "277662.6 (±48046.6)"
$ws.Range("M93").Value = "13.12 (±0.19)"
$ws.Range("N93").Value = "1170.4 (±202.5)"

# Row 94
$ws.Range("C94").Value = 29097.8
$ws.Range("E94").Value = 366059.2
$ws.Range("G94").Value = 12.58
$ws.Range("H94").Value = 0.14
$ws.Range("J94").Value = 1590.1
$ws.Range("L94").Value = "366059.2 (±54697.7)"
$ws.Range("M94").Value = "12.58 (±0.14)"
$ws.Range("N94").Value = "1590.1 (±237.6)"

# Row 95
$ws.Range("C95").Value = 50234
$ws.Range("D95").Value = 6800.1
$ws.Range("E95").Value = 659751.9999999999
$ws.Range("F95").Value = 99568.60000000001
$ws.Range("G95").Value = 13.13
$ws.Range("H95").Value = 0.18
$ws.Range("J95").Value = 1411.4
$ws.Range("K95").Value = 213
$ws.Range("L95").Value = "659752.0 (±99568.6)"
$ws.Range("M95").Value = "13.13 (±0.18)"
$ws.Range("N95").Value = "1411.4 (±213.0)"
